$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Friendly Matches"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2026-01-15"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "21:00:00"
$ws.Range("D2").Value = "Club Olimpia"
$ws.Range("E2").Value = "Colo Colo"
$ws.Range("F2").Value = 2.66
$ws.Range("G2").Value = 2.84
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 2.9
$ws.Range("K2").Value = 3.1
$ws.Range("L2").Value = 2.36
$ws.Range("M2").Value = 1.16
$ws.Range("N2").Value = 2.22
$ws.Range("O2").Value = 1.77
$ws.Range("P2").Value = 1.4
$ws.Range("Q2").Value = 3.25
$ws.Range("R2").Value = 1.13
$ws.Range("S2").Value = 7.8
$ws.Range("T2").Value = 2.26
$ws.Range("U2").Value = 1.57
$ws.Range("V2").Value = 1.4
$ws.Range("W2").Value = 1.56
$ws.Range("X2").Value = 8.800000000000001
$ws.Range("Y2").Value = 9.6
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 80
$ws.Range("AB2").Value = 7.8
$ws.Range("AC2").Value = 7.8
$ws.Range("AD2").Value = 18
$ws.Range("AE2").Value = 70
$ws.Range("AF2").Value = 17
$ws.Range("AG2").Value = 16
$ws.Range("AH2").Value = 36
$ws.Range("AI2").Value = 120
$ws.Range("AJ2").Value = 55
$ws.Range("AK2").Value = 55
$ws.Range("AL2").Value = 110
$ws.Range("AM2").Value = 390
$ws.Range("AN2").Value = 80
$ws.Range("AO2").Value = 100
